# Apply the edits described by the target diff:
#
# 1. The table on slide 16 gets a new table style (tableStyleId swapped
#    from {BBC196B3-D6D3-474B-B157-C7026572DCA4} to
#    {B4057EE5-0B90-4FD7-A865-190F19FABF8E}).
# 2. The presentation's theme colour palette changes from the "Integral"
#    scheme to the stock "Office Theme" scheme (the deck's live/applied
#    theme - ppt/theme/theme2.xml - is what the PowerPoint object model
#    exposes as Slide.ThemeColorScheme / NotesPage.ThemeColorScheme).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
        break
    }
}
if ($tableShape -ne $null) {
    $tableShape.Table.ApplyStyle("{B4057EE5-0B90-4FD7-A865-190F19FABF8E}")
}

# --- 2. Theme colours -------------------------------------------------
# Index order (matches MsoThemeColorSchemeIndex / <a:clrScheme> order):
#  1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3 8=accent4
#  9=accent5 10=accent6 11=hlink 12=folHlink
# RGB() packs as 0xBBGGRR, so the long value below is
# B*65536 + G*256 + R for each target "Office Theme" hex colour.
$officeThemeRGB = @(
    0x000000,  # dk1      000000
    0xFFFFFF,  # lt1      FFFFFF
    0x6A5444,  # dk2      44546A
    0xE6E6E7,  # lt2      E7E6E6
    0xD59B5B,  # accent1  5B9BD5
    0x317DED,  # accent2  ED7D31
    0xA5A5A5,  # accent3  A5A5A5
    0x00C0FF,  # accent4  FFC000
    0xC47244,  # accent5  4472C4
    0x47AD70,  # accent6  70AD47
    0xC16305,  # hlink    0563C1
    0x724F95   # folHlink 954F72
)

$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeRGB[$i - 1]
}
